# edit.ps1
# Applies the 2024-11-16 FlashScore "Jogos da Semana" workbook update:
#   - Odds refresh for the existing Blackpool - Northampton match (row 4)
#   - Odds refresh for the existing Castellon - Mirandes match (row 15)
#   - Two newly added SPAIN - LALIGA2 matches appended as rows 17 and 18
#     (Levante - Elche, Racing Santander - Burgos CF)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Blackpool - Northampton): odds refresh -------------------------
$ws.Range("G4").Value = 1.53
$ws.Range("H4").Value = 4.33
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.05
$ws.Range("K4").Value = 2.5
$ws.Range("S4").Value = 1.29
$ws.Range("T4").Value = 3.5
$ws.Range("Y4").Value = 8.5
$ws.Range("AC4").Value = 17
$ws.Range("AQ4").Value = 21
$ws.Range("AT4").Value = 3.5
$ws.Range("BB4").Value = 151

# --- Row 15 (Castellon - Mirandes): odds refresh ----------------------------
$ws.Range("Q15").Value = 1.88
$ws.Range("R15").Value = 1.98

# --- New rows 17 and 18: two additional SPAIN - LALIGA2 fixtures -----------
# Row 17: Levante vs Elche
$ws.Range("A17").Value = "jgo01QEs"
$ws.Range("B17").Value = "16/11/2024"
$ws.Range("C17").Value = "14:30"
$ws.Range("D17").Value = "SPAIN - LALIGA2"
$ws.Range("E17").Value = "Levante"
$ws.Range("F17").Value = "Elche"
$ws.Range("G17").Value = 2.38
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 3
$ws.Range("J17").Value = 3.1
$ws.Range("K17").Value = 2.05
$ws.Range("L17").Value = 3.75
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 9
$ws.Range("O17").Value = 1.36
$ws.Range("P17").Value = 3
$ws.Range("Q17").Value = 2.1
$ws.Range("R17").Value = 1.7
$ws.Range("S17").Value = 1.44
$ws.Range("T17").Value = 2.63
$ws.Range("U17").Value = 1.83
$ws.Range("V17").Value = 1.83
$ws.Range("W17").Value = 7.5
$ws.Range("X17").Value = 11
$ws.Range("Y17").Value = 10
$ws.Range("Z17").Value = 23
$ws.Range("AA17").Value = 21
$ws.Range("AB17").Value = 34
$ws.Range("AC17").Value = 8.5
$ws.Range("AD17").Value = 6
$ws.Range("AE17").Value = 15
$ws.Range("AF17").Value = 51
$ws.Range("AG17").Value = 301
$ws.Range("AH17").Value = 8.5
$ws.Range("AI17").Value = 15
$ws.Range("AJ17").Value = 11
$ws.Range("AK17").Value = 34
$ws.Range("AL17").Value = 26
$ws.Range("AM17").Value = 34
$ws.Range("AN17").Value = 4.33
$ws.Range("AO17").Value = 13
$ws.Range("AP17").Value = 26
$ws.Range("AQ17").Value = 41
$ws.Range("AR17").Value = 67
$ws.Range("AS17").Value = 201
$ws.Range("AT17").Value = 2.63
$ws.Range("AU17").Value = 8.5
$ws.Range("AV17").Value = 51
$ws.Range("AW17").Value = 5
$ws.Range("AX17").Value = 17
$ws.Range("AY17").Value = 29
$ws.Range("AZ17").Value = 51
$ws.Range("BA17").Value = 81
$ws.Range("BB17").Value = 201
$ws.Range("BC17").Value = 81
$ws.Range("BD17").Value = 81

# Row 18: Racing Santander vs Burgos CF
$ws.Range("A18").Value = "pfJSJYU8"
$ws.Range("B18").Value = "16/11/2024"
$ws.Range("C18").Value = "14:30"
$ws.Range("D18").Value = "SPAIN - LALIGA2"
$ws.Range("E18").Value = "Racing Santander"
$ws.Range("F18").Value = "Burgos CF"
$ws.Range("G18").Value = 1.55
$ws.Range("H18").Value = 4.2
$ws.Range("I18").Value = 5.5
$ws.Range("J18").Value = 2.1
$ws.Range("K18").Value = 2.38
$ws.Range("L18").Value = 5.5
$ws.Range("M18").Value = 1.04
$ws.Range("N18").Value = 13
$ws.Range("O18").Value = 1.22
$ws.Range("P18").Value = 4
$ws.Range("Q18").Value = 1.7
$ws.Range("R18").Value = 2.1
$ws.Range("S18").Value = 1.33
$ws.Range("T18").Value = 3.25
$ws.Range("U18").Value = 1.8
$ws.Range("V18").Value = 1.91
$ws.Range("W18").Value = 8
$ws.Range("X18").Value = 8
$ws.Range("Y18").Value = 8.5
$ws.Range("Z18").Value = 12
$ws.Range("AA18").Value = 12
$ws.Range("AB18").Value = 23
$ws.Range("AC18").Value = 13
$ws.Range("AD18").Value = 8
$ws.Range("AE18").Value = 17
$ws.Range("AF18").Value = 51
$ws.Range("AG18").Value = 201
$ws.Range("AH18").Value = 17
$ws.Range("AI18").Value = 29
$ws.Range("AJ18").Value = 17
$ws.Range("AK18").Value = 51
$ws.Range("AL18").Value = 41
$ws.Range("AM18").Value = 41
$ws.Range("AN18").Value = 3.6
$ws.Range("AO18").Value = 7.5
$ws.Range("AP18").Value = 17
$ws.Range("AQ18").Value = 23
$ws.Range("AR18").Value = 41
$ws.Range("AS18").Value = 126
$ws.Range("AT18").Value = 3.25
$ws.Range("AU18").Value = 8
$ws.Range("AV18").Value = 51
$ws.Range("AW18").Value = 7
$ws.Range("AX18").Value = 29
$ws.Range("AY18").Value = 34
$ws.Range("AZ18").Value = 101
$ws.Range("BA18").Value = 101
$ws.Range("BB18").Value = 201
$ws.Range("BC18").Value = 81
$ws.Range("BD18").Value = 81
